$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Update the quarterly rows. In the source data each year has four
# quarter rows (A, B, C, D). The B ("second quarter") and C ("third quarter")
# rows were recorded in the wrong order; swap their content (label + B/C/D/E
# values) so that B precedes C chronologically within every year block.
$ws.Range("A3").Value = "2000年C"
$ws.Range("B3").Value = 99.7
$ws.Range("D3").Value = 12.7
$ws.Range("E3").Value = 4417
$ws.Range("A4").Value = "2000年B"
$ws.Range("B4").Value = 100.1
$ws.Range("D4").Value = -1.8
$ws.Range("E4").Value = 2738.5

$ws.Range("A7").Value = "2001年C"
$ws.Range("B7").Value = 99.7
$ws.Range("D7").Value = 18.7
$ws.Range("E7").Value = 4077.4
$ws.Range("A8").Value = "2001年B"
$ws.Range("B8").Value = 96.90000000000001
$ws.Range("D8").Value = 87.7
$ws.Range("E8").Value = 2636.4

$ws.Range("A11").Value = "2002年C"
$ws.Range("B11").Value = 99.59999999999999
$ws.Range("D11").Value = 16.1
$ws.Range("E11").Value = 3341.9
$ws.Range("A12").Value = "2002年B"
$ws.Range("B12").Value = 99.59999999999999
$ws.Range("D12").Value = 9.5
$ws.Range("E12").Value = 2087.5

$ws.Range("A15").Value = "2003年C"
$ws.Range("B15").Value = 100.1
$ws.Range("C15").Value = 0.5
$ws.Range("D15").Value = -3.3
$ws.Range("E15").Value = 4167.6
$ws.Range("A16").Value = "2003年B"
$ws.Range("B16").Value = 101.2
$ws.Range("C16").Value = 1.6
$ws.Range("D16").Value = -11.1
$ws.Range("E16").Value = 2683.8

$ws.Range("A19").Value = "2004年C"
$ws.Range("B19").Value = 98.90000000000001
$ws.Range("C19").Value = -1.2
$ws.Range("D19").Value = 46.1
$ws.Range("E19").Value = 4657.3
$ws.Range("A20").Value = "2004年B"
$ws.Range("B20").Value = 99.2
$ws.Range("C20").Value = -2
$ws.Range("D20").Value = 21.7
$ws.Range("E20").Value = 2849.6

$ws.Range("A23").Value = "2005年C"
$ws.Range("B23").Value = 100.4
$ws.Range("C23").Value = 1.5
$ws.Range("D23").Value = -16.8
$ws.Range("E23").Value = 4880.5
$ws.Range("A24").Value = "2005年B"
$ws.Range("B24").Value = 99.90000000000001
$ws.Range("C24").Value = 0.7
$ws.Range("D24").Value = 1.7
$ws.Range("E24").Value = 3971.7

$ws.Range("A27").Value = "2006年C"
$ws.Range("B27").Value = 98.3
$ws.Range("C27").Value = -2.1
$ws.Range("D27").Value = 123.9
$ws.Range("E27").Value = 5836.2
$ws.Range("A28").Value = "2006年B"
$ws.Range("B28").Value = 99.7
$ws.Range("C28").Value = -0.2
$ws.Range("D28").Value = 15.4
$ws.Range("E28").Value = 3821.5

$ws.Range("A31").Value = "2007年C"
$ws.Range("B31").Value = 101.5
$ws.Range("C31").Value = 3.2
$ws.Range("D31").Value = -31
$ws.Range("E31").Value = 5825
$ws.Range("A32").Value = "2007年B"
$ws.Range("B32").Value = 103.4
$ws.Range("C32").Value = 3.7
$ws.Range("D32").Value = -41.8
$ws.Range("E32").Value = 3570.2

$ws.Range("A35").Value = "2008年C"
$ws.Range("B35").Value = 97.5
$ws.Range("C35").Value = -4.7
$ws.Range("D35").Value = 106.8
$ws.Range("E35").Value = 5911.9
$ws.Range("A36").Value = "2008年B"
$ws.Range("B36").Value = 96.59999999999999
$ws.Range("C36").Value = -6.3
$ws.Range("D36").Value = 86.59999999999999
$ws.Range("E36").Value = 3767.6

$ws.Range("A39").Value = "2009年C"
$ws.Range("B39").Value = 99.90000000000001
$ws.Range("C39").Value = 2.4
$ws.Range("D39").Value = 2.7
$ws.Range("E39").Value = 5639.9
$ws.Range("A40").Value = "2009年B"
$ws.Range("B40").Value = 99.7
$ws.Range("C40").Value = 2.8
$ws.Range("D40").Value = 4
$ws.Range("E40").Value = 3348.2

$ws.Range("A43").Value = "2010年C"
$ws.Range("B43").Value = 100.2
$ws.Range("C43").Value = 1.5
$ws.Range("D43").Value = -4.5
$ws.Range("E43").Value = 6499.6
$ws.Range("A44").Value = "2010年B"
$ws.Range("B44").Value = 100.7
$ws.Range("C44").Value = 1.2
$ws.Range("D44").Value = -9.699999999999999
$ws.Range("E44").Value = 4090.4

$ws.Range("A47").Value = "2011年C"
$ws.Range("B47").Value = 103.3
$ws.Range("C47").Value = 2.8
$ws.Range("D47").Value = -43.6
$ws.Range("E47").Value = 6401.7
$ws.Range("A48").Value = "2011年B"
$ws.Range("B48").Value = 102.5
$ws.Range("C48").Value = 1.8
$ws.Range("D48").Value = -20.6
$ws.Range("E48").Value = 4025.2

$ws.Range("A51").Value = "2012年C"
$ws.Range("B51").Value = 83
$ws.Range("C51").Value = -11.5
$ws.Range("D51").Value = 11.4
$ws.Range("E51").Value = 5641.1
$ws.Range("A52").Value = "2012年B"
$ws.Range("B52").Value = 99.59999999999999
$ws.Range("C52").Value = -3.3
$ws.Range("D52").Value = 11.4
$ws.Range("E52").Value = 3839

$ws.Range("A55").Value = "2013年C"
$ws.Range("B55").Value = 98.59999999999999
$ws.Range("C55").Value = -2.8
$ws.Range("D55").Value = 40.6
$ws.Range("E55").Value = 4018.3
$ws.Range("A56").Value = "2013年B"
$ws.Range("B56").Value = 100.4
$ws.Range("C56").Value = 1.6
$ws.Range("D56").Value = -19.4
$ws.Range("E56").Value = 2708.6

$ws.Range("A59").Value = "2014年C"
$ws.Range("B59").Value = 103.1
$ws.Range("C59").Value = 3.8
$ws.Range("D59").Value = -36.4
$ws.Range("E59").Value = 2384.3
$ws.Range("A60").Value = "2014年B"
$ws.Range("B60").Value = 102.8
$ws.Range("C60").Value = 0.4
$ws.Range("D60").Value = -21.6
$ws.Range("E60").Value = 1564.6

$ws.Range("A63").Value = "2015年C"
$ws.Range("B63").Value = 99.40000000000001
$ws.Range("C63").Value = -3.7
$ws.Range("D63").Value = 12.4
$ws.Range("E63").Value = 2127.2
$ws.Range("A64").Value = "2015年B"
$ws.Range("B64").Value = 98.2
$ws.Range("C64").Value = -4.6
$ws.Range("D64").Value = 24.2
$ws.Range("E64").Value = 1348.4

# Step 2: The workbook previously duplicated two columns ("F" = 产销率,
# "G" = 销售量) that just mirrored columns B and E. Remove them entirely;
# deleting (rather than merely clearing) shifts nothing else and lets Excel
# recompute the sheet dimension automatically (becomes A1:E65).
$ws.Range("F1:G65").Delete()
